$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the quality/star values (balance card attr of level, star and quality)
$ws.Range("F4").Value = 3
$ws.Range("F6").Value = 3

# Update the active selection to F5
$ws.Range("F5").Select()
